# ADD results from server
# Update computed result values (row 2) on each yearly sheet with new
# server-provided figures.

$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param(
        [string]$SheetName,
        [hashtable]$Cells
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($cellRef in $Cells.Keys) {
        $ws.Range($cellRef).Value = $Cells[$cellRef]
    }
}

Set-RowValues "2025" @{
    "B2" = 1037.265132737054
    "E2" = 28926.05393052954
    "G2" = 8095.925712661834
    "I2" = 16171.06685703679
    "L2" = 48492.22142001599
    "M2" = 10595.37713982
    "N2" = 7071.74531360843
    "O2" = 6993.890772562212
}

Set-RowValues "2030" @{
    "A2" = 0
    "B2" = 4157.588990853394
    "E2" = 45991.90904307188
    "G2" = 8095.925712661834
    "I2" = 37079.12819938764
    "L2" = 54844.03303316472
    "M2" = 17449.04999683176
    "N2" = 9024.733389685653
    "O2" = 9724.258249348202
}

Set-RowValues "2035" @{
    "A2" = 2754.31755456332
    "B2" = 6368.910634126893
    "E2" = 57457.45307013817
    "G2" = 8095.925712661834
    "I2" = 52465.73681402855
    "L2" = 54844.03303316472
    "M2" = 21912.87293902603
    "N2" = 13034.3101291405
    "O2" = 12860.17168993684
}

Set-RowValues "2040" @{
    "A2" = 2754.31755456332
    "B2" = 6368.910634126893
    "E2" = 57457.45307013817
    "G2" = 8095.925712661834
    "I2" = 52465.73681402855
    "L2" = 54844.03303316472
    "M2" = 21912.87293902603
    "N2" = 13151.8694977663
    "O2" = 12860.17168993684
}

Set-RowValues "2045" @{
    "A2" = 5713.151062849596
    "B2" = 6368.910634126893
    "E2" = 57457.45307013817
    "G2" = 8095.925712661834
    "I2" = 52465.73681402855
    "L2" = 54844.03303316472
    "M2" = 21912.87293902603
    "N2" = 13601.08685191924
    "O2" = 14937.1305943757
}

Set-RowValues "2050" @{
    "A2" = 5713.151062849596
    "B2" = 6368.910634126893
    "E2" = 57457.45307013817
    "G2" = 8095.925712661834
    "I2" = 52465.73681402855
    "L2" = 54844.03303316472
    "M2" = 21912.87293902603
    "N2" = 13601.08685191924
    "O2" = 14937.1305943757
}
